$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.176.90"
$ws.Range("E2").Value = "  +13.02%  "

$ws.Range("D3").Value = "'1.779.17"
$ws.Range("E3").Value = "  +6.44%  "

$ws.Range("D4").Value = "'0.992"
$ws.Range("E4").Value = "  -0.54%  "

$ws.Range("D5").Value = "'230.27"
$ws.Range("E5").Value = "  +5.23%  "

$ws.Range("D6").Value = "'0.550"
$ws.Range("E6").Value = "  +5.53%  "

$ws.Range("D7").Value = "'0.993"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "'31.57"
$ws.Range("E8").Value = "  +6.16%  "

$ws.Range("D9").Value = "'46.26"
$ws.Range("E9").Value = "  +5.03%  "

$ws.Range("E10").Value = "  +5.37%  "

$ws.Range("D11").Value = "'0.0670"
$ws.Range("E11").Value = "  +8.61%  "

$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").Value = "'2.027.49"
$ws.Range("E13").Value = "  +6.02%  "

$ws.Range("D14").Value = "'1.773.47"
$ws.Range("E14").Value = "  +6.05%  "

$ws.Range("D15").Value = "'0.637"
$ws.Range("E15").Value = "  +3.71%  "

$ws.Range("D16").Value = "'10.40"
$ws.Range("E16").Value = "  -3.87%  "

$ws.Range("D17").Value = "'33.991.79"
$ws.Range("E17").Value = "  +12.32%  "

$ws.Range("D18").Value = "'4.32"
$ws.Range("E18").Value = "  +7.99%  "

$ws.Range("D19").Value = "'69.75"
$ws.Range("E19").Value = "  +6.38%  "

$ws.Range("D20").Value = "'263.61"
$ws.Range("E20").Value = "  +6.43%  "

$ws.Range("D21").Value = "'0.0₃0756"
$ws.Range("E21").Value = "  +5.62%  "

$ws.Range("D22").Value = "'0.995"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").Value = "'10.45"
$ws.Range("E23").Value = "  +4.07%  "

$ws.Range("D24").Value = "'4.40"
$ws.Range("E24").Value = "  +1.72%  "

$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "'161.38"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").Value = "'16.72"
$ws.Range("E27").Value = "  +5.33%  "

$ws.Range("D28").Value = "'0.116"
$ws.Range("E28").Value = "  +4.39%  "

$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  +5.49%  "

$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").Value = "'3.82"
$ws.Range("E31").Value = "  +10.32%  "

$ws.Range("D32").Value = "'0.0512"
$ws.Range("E32").Value = "  +2.75%  "

$ws.Range("E33").Value = "  +5.46%  "

$ws.Range("D34").Value = "'3.56"
$ws.Range("E34").Value = "  +8.76%  "

$ws.Range("D35").Value = "'1.568.56"
$ws.Range("E35").Value = "  +6.72%  "

$ws.Range("E36").Value = "  +5.38%  "

$ws.Range("D37").Value = "'88.46"
$ws.Range("E37").Value = "  +9.99%  "

$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  +2.43%  "

$ws.Range("D39").Value = "'0.627"
$ws.Range("E39").Value = "  +6.96%  "

$ws.Range("D40").Value = "'0.0187"
$ws.Range("E40").Value = "  +4.96%  "

$ws.Range("D41").Value = "'2.79"
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("E42").Value = "  +2.24%  "

$ws.Range("D43").Value = "'0.910"
$ws.Range("E43").Value = "  +6.25%  "

$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  +4.67%  "

$ws.Range("D45").Value = "'0.0517"
$ws.Range("E45").Value = "  +3.00%  "

$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'54.85"
$ws.Range("E46").Value = "  +3.74%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.05"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'1.918.83"
$ws.Range("E48").Value = "  +6.00%  "

$ws.Range("D49").Value = "'5.74"
$ws.Range("E49").Value = "  +5.90%  "

$ws.Range("D50").Value = "'0.996"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'0.0₆0122"
$ws.Range("E51").Value = "  +3.89%  "
